# no-op for now
$p = $ppt.ActivePresentation
Write-Host "slides count: $($p.Slides.Count)"
